$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts the existing rows 13..88
# down to 14..89 (and carries the row 12 formatting down onto the new row,
# matching the date-styled column D).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new data record.
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44537
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112022
$ws.Range("G13").Value = "Arveja Verde"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 61
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 15984
$ws.Range("N13").Value = '$/saco 25 kilos'
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 639
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
